# Actualización desde MV -datos-
# Adds a new survey-vintage column (BH, "Agosto.2021") to the table and a
# new observation row (75, period "01-04-2021"), and records a revision to
# the "01-01-2021" row's latest estimate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BH (column 60): header + copied-forward values ---------

# Copy BG1's formatting (bold, centered, bordered header style) into BH1,
# then overwrite its value with the new vintage label.
$ws.Cells.Item(1, 59).Copy($ws.Cells.Item(1, 60))
$ws.Cells.Item(1, 60).Value2 = "Agosto.2021"

# Rows 2-73: the new vintage simply repeats the latest known estimate
# (column BG) for every already-closed reference period.
for ($r = 2; $r -le 73; $r++) {
    $ws.Cells.Item($r, 60).Value2 = $ws.Cells.Item($r, 59).Value2
}

# Row 74 ("01-01-2021"): the new vintage brings a revised estimate.
$ws.Cells.Item(74, 60).Value2 = 34660

# --- New row 75: newest reference period, only the new vintage filled --

# Force text (not an auto-converted date serial) for the period label, then
# drop back to the default style so the cell matches its unstyled peers in
# column A.
$ws.Cells.Item(75, 1).NumberFormat = "@"
$ws.Cells.Item(75, 1).Value2 = "01-04-2021"
$ws.Cells.Item(75, 1).Style = "Normal"

$ws.Cells.Item(75, 60).Value2 = 36066
